# "Generate Report for Archive"
# The localization status report was regenerated: the "Ready for handoff"
# status became "In Translation" everywhere it appears, and the now-narrower
# "Status" columns were shrunk to fit the shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns are E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:E4").Value = "In Translation"
$overview.Range("F2:F4").Value = "In Translation"

# --- zh-cn / de-de detail sheets: status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# Shrink the status columns to fit the shorter "In Translation" text
# (was width 17.2159881591797, now ~13.41).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
